# 245-MS-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME-Loanproduct.xlsx
# "multi browser implementation for chrome and firefox and accounting cash"
#
# The loan-product input sheet gains a new "currency" / "US Dollar" row
# (row 6), that row's formatting is normalised to the plain green
# "value" look used by the rest of the sheet, the now entirely-empty
# helper column C is removed, and the input sheet (rather than the
# output sheet) becomes the active/selected tab & range.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# 1. Populate the new currency row.
$ws1.Range("A6").Value = "currency"
$ws1.Range("B6").Value = "US Dollar"

# 2. Re-style B6 from the bold/grey "pair" look it inherited from the old
#    Currency/US Dollar header cells to the plain green look used by the
#    other single-value cells (e.g. B10) elsewhere on the sheet.
$ws1.Range("B10").Copy()
$ws1.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Column C only ever held empty placeholder cells (C5, C6, C12) - drop it.
$ws1.Range("C:C").Delete()

# 4. Switch the active tab/selection from the output sheet to the input
#    sheet, highlighting the newly added currency row.
$ws1.Activate()
$ws1.Range("A6:B6").Select()
